$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Avij input variables
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 14

# Update the H21 formula (B2*B5*B5 -> B2*B5*(B5-1))
$ws.Range("H21").Formula = "=B2*B5*(B5-1)"

# Move the active selection to H21
$ws.Range("H21").Select()
